$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 50
$ws.Range("B6").Value = "Solomito_Carne"
$ws.Range("C6").Value = 20000
$ws.Range("D6").Value = 32
